# TC05-04-053 AddMyTicket - add new worksheet with ticket test data
$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the workbook (after the last existing sheet)
$sheetCount = $wb.Worksheets.Count
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($sheetCount))
$newSheet.Name = "TC05-04-053 AddMyTicket"

# Header row
$newSheet.Range("A1").Value = "ticketType"
$newSheet.Range("B1").Value = "category"
$newSheet.Range("C1").Value = "priority"
$newSheet.Range("D1").Value = "expectedDate"
$newSheet.Range("E1").Value = "subject"
$newSheet.Range("F1").Value = "description"
$newSheet.Range("G1").Value = "document"

# Row 2
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 2
$newSheet.Range("C2").Value = 3
$newSheet.Range("D2").Value = "'05/02/2024"
$newSheet.Range("E2").Value = "Ticket 1saa"
$newSheet.Range("F2").Value = "issa ticket"
$newSheet.Range("G2").Value = "C:\Users\Hasnul\Katalon Studio\Katalon-Studio-Training-Assessment\Empxtrack_Nafis\File Upload Test Data\Katalon Studio Training Day 2.pdf"

# Row 3
$newSheet.Range("A3").Value = 2
$newSheet.Range("E3").Value = "Ticket 2aa"
$newSheet.Range("F3").Value = "issa ticketttttt"
$newSheet.Range("G3").Value = "C:\Users\Hasnul\Katalon Studio\Katalon-Studio-Training-Assessment\Empxtrack_Nafis\File Upload Test Data\Katalon Studio Training Day 2.pdf"

# Make the new sheet active, with H5 selected (matches the new sheet's saved view)
$newSheet.Activate()
$newSheet.Range("H5").Select()

Write-Host "done"
